$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was old row 10's data for columns D, M, N, O, P, Q, R, S, T)
$ws.Range("D2").Value = 44174
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 3700
$ws.Range("O2").Value = 3800
$ws.Range("P2").Value = 3747
$ws.Range("Q2").Value = '$/bandeja 2 kilos'
$ws.Range("R2").Value = "Provincia de Linares"
$ws.Range("S2").Value = 1874
$ws.Range("T2").Value = 2

# Row 3 (was old row 7's data)
$ws.Range("D3").Value = 44540
$ws.Range("M3").Value = 240
$ws.Range("N3").Value = 3500
$ws.Range("O3").Value = 3800
$ws.Range("P3").Value = 3650
$ws.Range("Q3").Value = '$/bandeja 2 kilos'
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 1825
$ws.Range("T3").Value = 2

# Row 5 (was old row 2's data)
$ws.Range("D5").Value = 44181
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 3600
$ws.Range("O5").Value = 3800
$ws.Range("P5").Value = 3692
$ws.Range("Q5").Value = '$/bandeja 2 kilos'
$ws.Range("R5").Value = "Provincia de Diguillín"
$ws.Range("S5").Value = 1846
$ws.Range("T5").Value = 2

# Row 6 (was old row 3's data)
$ws.Range("D6").Value = 44181
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 1800
$ws.Range("O6").Value = 2000
$ws.Range("P6").Value = 1875
$ws.Range("Q6").Value = '$/envase 1 kilo'
$ws.Range("R6").Value = "Provincia de Diguillín"
$ws.Range("S6").Value = 1875
$ws.Range("T6").Value = 1

# Row 7 (was old row 6's data)
$ws.Range("D7").Value = 44539
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 3800
$ws.Range("O7").Value = 4000
$ws.Range("P7").Value = 3900
$ws.Range("Q7").Value = '$/bandeja 2 kilos'
$ws.Range("R7").Value = "Región del Maule"
$ws.Range("S7").Value = 1950
$ws.Range("T7").Value = 2

# Row 10 (was old row 5's data)
$ws.Range("D10").Value = 44594
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 2800
$ws.Range("P10").Value = 2650
$ws.Range("Q10").Value = '$/bandeja 2 kilos'
$ws.Range("R10").Value = "Provincia de Linares"
$ws.Range("S10").Value = 1325
$ws.Range("T10").Value = 2
